# Commit "updated order files with png": the stimulus image filenames referenced
# in columns B (Condition/CUE marker), D (Filename_Left) and E (Filename_Right)
# were renamed from *.jpg to *.png. Apply the same renames to every affected
# cell, then restore the worksheet selection to the full used range.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'CUE'
$ws.Range('D2').Value = 'CoverCue.png'
$ws.Range('E2').Value = 'CoverCue.png'
$ws.Range('D4').Value = 'Face9_L.png'
$ws.Range('E4').Value = 'Face9_L.png'
$ws.Range('D5').Value = 'Face6_L.png'
$ws.Range('E5').Value = 'Face6_L.png'
$ws.Range('D6').Value = 'Face18_L.png'
$ws.Range('E6').Value = 'Face18_L.png'
$ws.Range('D7').Value = 'Face13_L.png'
$ws.Range('E7').Value = 'Face13_L.png'
$ws.Range('D8').Value = 'Face12_L.png'
$ws.Range('E8').Value = 'Face12_L.png'
$ws.Range('D9').Value = 'Face3_L.png'
$ws.Range('E9').Value = 'Face3_L.png'
$ws.Range('D10').Value = 'Face22_L.png'
$ws.Range('E10').Value = 'Face22_L.png'
$ws.Range('D11').Value = 'Face21_L.png'
$ws.Range('E11').Value = 'Face21_L.png'
$ws.Range('D12').Value = 'Face11_L.png'
$ws.Range('E12').Value = 'Face11_L.png'
$ws.Range('D13').Value = 'Face8_L.png'
$ws.Range('E13').Value = 'Face8_L.png'
$ws.Range('D14').Value = 'Face5_L.png'
$ws.Range('E14').Value = 'Face5_L.png'
$ws.Range('D15').Value = 'Face5_L.png'
$ws.Range('E15').Value = 'Face5_L.png'
$ws.Range('D16').Value = 'Face15_L.png'
$ws.Range('E16').Value = 'Face15_L.png'
$ws.Range('D17').Value = 'Face10_L.png'
$ws.Range('E17').Value = 'Face10_L.png'
$ws.Range('D18').Value = 'Face23_L.png'
$ws.Range('E18').Value = 'Face23_L.png'
$ws.Range('D19').Value = 'Face16_L.png'
$ws.Range('E19').Value = 'Face16_L.png'
$ws.Range('D20').Value = 'Face20_L.png'
$ws.Range('E20').Value = 'Face20_L.png'
$ws.Range('D21').Value = 'Face7_L.png'
$ws.Range('E21').Value = 'Face7_L.png'
$ws.Range('D22').Value = 'Face4_L.png'
$ws.Range('E22').Value = 'Face4_L.png'
$ws.Range('D23').Value = 'Face14_L.png'
$ws.Range('E23').Value = 'Face14_L.png'
$ws.Range('D24').Value = 'Face19_L.png'
$ws.Range('E24').Value = 'Face19_L.png'
$ws.Range('D25').Value = 'Face2_L.png'
$ws.Range('E25').Value = 'Face2_L.png'
$ws.Range('D26').Value = 'Face17_L.png'
$ws.Range('E26').Value = 'Face17_L.png'
$ws.Range('D27').Value = 'Face1_L.png'
$ws.Range('E27').Value = 'Face1_L.png'
$ws.Range('B28').Value = 'CUE'
$ws.Range('D28').Value = 'UncoverCue.png'
$ws.Range('E28').Value = 'UncoverCue.png'
$ws.Range('B30').Value = 'Pseudo'
$ws.Range('D30').Value = 'Face17_R.png'
$ws.Range('E30').Value = 'Face17_L.png'
$ws.Range('B31').Value = 'Pseudo'
$ws.Range('D31').Value = 'Face9_R.png'
$ws.Range('E31').Value = 'Face9_L.png'
$ws.Range('B32').Value = 'Pseudo'
$ws.Range('D32').Value = 'Face20_R.png'
$ws.Range('E32').Value = 'Face20_L.png'
$ws.Range('B33').Value = 'Pseudo'
$ws.Range('D33').Value = 'Face8_R.png'
$ws.Range('E33').Value = 'Face8_L.png'
$ws.Range('B34').Value = 'Pseudo'
$ws.Range('D34').Value = 'Face23_R.png'
$ws.Range('E34').Value = 'Face23_L.png'
$ws.Range('B35').Value = 'Pseudo'
$ws.Range('D35').Value = 'Face11_R.png'
$ws.Range('E35').Value = 'Face11_L.png'
$ws.Range('B36').Value = 'Pseudo'
$ws.Range('D36').Value = 'Face1_R.png'
$ws.Range('E36').Value = 'Face1_L.png'
$ws.Range('B37').Value = 'Pseudo'
$ws.Range('D37').Value = 'Face12_R.png'
$ws.Range('E37').Value = 'Face12_L.png'
$ws.Range('B38').Value = 'Pseudo'
$ws.Range('D38').Value = 'Face22_R.png'
$ws.Range('E38').Value = 'Face22_L.png'
$ws.Range('B39').Value = 'Pseudo'
$ws.Range('D39').Value = 'Face18_R.png'
$ws.Range('E39').Value = 'Face18_L.png'
$ws.Range('B40').Value = 'Pseudo'
$ws.Range('D40').Value = 'Face13_R.png'
$ws.Range('E40').Value = 'Face13_L.png'
$ws.Range('B41').Value = 'Pseudo'
$ws.Range('D41').Value = 'Face2_R.png'
$ws.Range('E41').Value = 'Face2_L.png'
$ws.Range('B42').Value = 'Pseudo'
$ws.Range('D42').Value = 'Face16_R.png'
$ws.Range('E42').Value = 'Face16_L.png'
$ws.Range('B43').Value = 'Pseudo'
$ws.Range('D43').Value = 'Face10_R.png'
$ws.Range('E43').Value = 'Face10_L.png'
$ws.Range('B44').Value = 'Pseudo'
$ws.Range('D44').Value = 'Face21_R.png'
$ws.Range('E44').Value = 'Face21_L.png'
$ws.Range('B45').Value = 'Pseudo'
$ws.Range('D45').Value = 'Face5_R.png'
$ws.Range('E45').Value = 'Face5_L.png'
$ws.Range('B46').Value = 'Pseudo'
$ws.Range('D46').Value = 'Face7_R.png'
$ws.Range('E46').Value = 'Face7_L.png'
$ws.Range('B47').Value = 'Pseudo'
$ws.Range('D47').Value = 'Face7_R.png'
$ws.Range('E47').Value = 'Face7_L.png'
$ws.Range('B48').Value = 'Pseudo'
$ws.Range('D48').Value = 'Face19_R.png'
$ws.Range('E48').Value = 'Face19_L.png'
$ws.Range('B49').Value = 'Pseudo'
$ws.Range('D49').Value = 'Face15_R.png'
$ws.Range('E49').Value = 'Face15_L.png'
$ws.Range('B50').Value = 'Pseudo'
$ws.Range('D50').Value = 'Face4_R.png'
$ws.Range('E50').Value = 'Face4_L.png'
$ws.Range('B51').Value = 'Pseudo'
$ws.Range('D51').Value = 'Face14_R.png'
$ws.Range('E51').Value = 'Face14_L.png'
$ws.Range('B52').Value = 'Pseudo'
$ws.Range('D52').Value = 'Face6_R.png'
$ws.Range('E52').Value = 'Face6_L.png'
$ws.Range('B53').Value = 'Pseudo'
$ws.Range('D53').Value = 'Face3_R.png'
$ws.Range('E53').Value = 'Face3_L.png'
$ws.Range('B55').Value = '3D'
$ws.Range('D55').Value = 'Face17_L.png'
$ws.Range('E55').Value = 'Face17_R.png'
$ws.Range('B56').Value = '3D'
$ws.Range('D56').Value = 'Face17_L.png'
$ws.Range('E56').Value = 'Face17_R.png'
$ws.Range('B57').Value = '3D'
$ws.Range('D57').Value = 'Face16_L.png'
$ws.Range('E57').Value = 'Face16_R.png'
$ws.Range('B58').Value = '3D'
$ws.Range('D58').Value = 'Face13_L.png'
$ws.Range('E58').Value = 'Face13_R.png'
$ws.Range('B59').Value = '3D'
$ws.Range('D59').Value = 'Face2_L.png'
$ws.Range('E59').Value = 'Face2_R.png'
$ws.Range('B60').Value = '3D'
$ws.Range('D60').Value = 'Face21_L.png'
$ws.Range('E60').Value = 'Face21_R.png'
$ws.Range('B61').Value = '3D'
$ws.Range('D61').Value = 'Face5_L.png'
$ws.Range('E61').Value = 'Face5_R.png'
$ws.Range('B62').Value = '3D'
$ws.Range('D62').Value = 'Face6_L.png'
$ws.Range('E62').Value = 'Face6_R.png'
$ws.Range('B63').Value = '3D'
$ws.Range('D63').Value = 'Face4_L.png'
$ws.Range('E63').Value = 'Face4_R.png'
$ws.Range('B64').Value = '3D'
$ws.Range('D64').Value = 'Face20_L.png'
$ws.Range('E64').Value = 'Face20_R.png'
$ws.Range('B65').Value = '3D'
$ws.Range('D65').Value = 'Face3_L.png'
$ws.Range('E65').Value = 'Face3_R.png'
$ws.Range('B66').Value = '3D'
$ws.Range('D66').Value = 'Face12_L.png'
$ws.Range('E66').Value = 'Face12_R.png'
$ws.Range('B67').Value = '3D'
$ws.Range('D67').Value = 'Face7_L.png'
$ws.Range('E67').Value = 'Face7_R.png'
$ws.Range('B68').Value = '3D'
$ws.Range('D68').Value = 'Face23_L.png'
$ws.Range('E68').Value = 'Face23_R.png'
$ws.Range('B69').Value = '3D'
$ws.Range('D69').Value = 'Face9_L.png'
$ws.Range('E69').Value = 'Face9_R.png'
$ws.Range('B70').Value = '3D'
$ws.Range('D70').Value = 'Face10_L.png'
$ws.Range('E70').Value = 'Face10_R.png'
$ws.Range('B71').Value = '3D'
$ws.Range('D71').Value = 'Face14_L.png'
$ws.Range('E71').Value = 'Face14_R.png'
$ws.Range('B72').Value = '3D'
$ws.Range('D72').Value = 'Face1_L.png'
$ws.Range('E72').Value = 'Face1_R.png'
$ws.Range('B73').Value = '3D'
$ws.Range('D73').Value = 'Face8_L.png'
$ws.Range('E73').Value = 'Face8_R.png'
$ws.Range('B74').Value = '3D'
$ws.Range('D74').Value = 'Face19_L.png'
$ws.Range('E74').Value = 'Face19_R.png'
$ws.Range('B75').Value = '3D'
$ws.Range('D75').Value = 'Face22_L.png'
$ws.Range('E75').Value = 'Face22_R.png'
$ws.Range('B76').Value = '3D'
$ws.Range('D76').Value = 'Face11_L.png'
$ws.Range('E76').Value = 'Face11_R.png'
$ws.Range('B77').Value = '3D'
$ws.Range('D77').Value = 'Face18_L.png'
$ws.Range('E77').Value = 'Face18_R.png'
$ws.Range('B78').Value = '3D'
$ws.Range('D78').Value = 'Face15_L.png'
$ws.Range('E78').Value = 'Face15_R.png'
$ws.Range('B80').Value = '2D'
$ws.Range('D80').Value = 'Face4_L.png'
$ws.Range('E80').Value = 'Face4_L.png'
$ws.Range('B81').Value = '2D'
$ws.Range('D81').Value = 'Face20_L.png'
$ws.Range('E81').Value = 'Face20_L.png'
$ws.Range('B82').Value = '2D'
$ws.Range('D82').Value = 'Face12_L.png'
$ws.Range('E82').Value = 'Face12_L.png'
$ws.Range('B83').Value = '2D'
$ws.Range('D83').Value = 'Face10_L.png'
$ws.Range('E83').Value = 'Face10_L.png'
$ws.Range('B84').Value = '2D'
$ws.Range('D84').Value = 'Face13_L.png'
$ws.Range('E84').Value = 'Face13_L.png'
$ws.Range('B85').Value = '2D'
$ws.Range('D85').Value = 'Face8_L.png'
$ws.Range('E85').Value = 'Face8_L.png'
$ws.Range('B86').Value = '2D'
$ws.Range('D86').Value = 'Face9_L.png'
$ws.Range('E86').Value = 'Face9_L.png'
$ws.Range('B87').Value = '2D'
$ws.Range('D87').Value = 'Face3_L.png'
$ws.Range('E87').Value = 'Face3_L.png'
$ws.Range('B88').Value = '2D'
$ws.Range('D88').Value = 'Face19_L.png'
$ws.Range('E88').Value = 'Face19_L.png'
$ws.Range('B89').Value = '2D'
$ws.Range('D89').Value = 'Face5_L.png'
$ws.Range('E89').Value = 'Face5_L.png'
$ws.Range('B90').Value = '2D'
$ws.Range('D90').Value = 'Face21_L.png'
$ws.Range('E90').Value = 'Face21_L.png'
$ws.Range('B91').Value = '2D'
$ws.Range('D91').Value = 'Face14_L.png'
$ws.Range('E91').Value = 'Face14_L.png'
$ws.Range('B92').Value = '2D'
$ws.Range('D92').Value = 'Face17_L.png'
$ws.Range('E92').Value = 'Face17_L.png'
$ws.Range('B93').Value = '2D'
$ws.Range('D93').Value = 'Face11_L.png'
$ws.Range('E93').Value = 'Face11_L.png'
$ws.Range('B94').Value = '2D'
$ws.Range('D94').Value = 'Face22_L.png'
$ws.Range('E94').Value = 'Face22_L.png'
$ws.Range('B95').Value = '2D'
$ws.Range('D95').Value = 'Face6_L.png'
$ws.Range('E95').Value = 'Face6_L.png'
$ws.Range('B96').Value = '2D'
$ws.Range('D96').Value = 'Face1_L.png'
$ws.Range('E96').Value = 'Face1_L.png'
$ws.Range('B97').Value = '2D'
$ws.Range('D97').Value = 'Face1_L.png'
$ws.Range('E97').Value = 'Face1_L.png'
$ws.Range('B98').Value = '2D'
$ws.Range('D98').Value = 'Face16_L.png'
$ws.Range('E98').Value = 'Face16_L.png'
$ws.Range('B99').Value = '2D'
$ws.Range('D99').Value = 'Face7_L.png'
$ws.Range('E99').Value = 'Face7_L.png'
$ws.Range('B100').Value = '2D'
$ws.Range('D100').Value = 'Face2_L.png'
$ws.Range('E100').Value = 'Face2_L.png'
$ws.Range('B101').Value = '2D'
$ws.Range('D101').Value = 'Face15_L.png'
$ws.Range('E101').Value = 'Face15_L.png'
$ws.Range('B102').Value = '2D'
$ws.Range('D102').Value = 'Face23_L.png'
$ws.Range('E102').Value = 'Face23_L.png'
$ws.Range('B103').Value = '2D'
$ws.Range('D103').Value = 'Face18_L.png'
$ws.Range('E103').Value = 'Face18_L.png'
$ws.Range('B105').Value = '3D'
$ws.Range('D105').Value = 'Face4_L.png'
$ws.Range('E105').Value = 'Face4_R.png'
$ws.Range('B106').Value = '3D'
$ws.Range('D106').Value = 'Face15_L.png'
$ws.Range('E106').Value = 'Face15_R.png'
$ws.Range('B107').Value = '3D'
$ws.Range('D107').Value = 'Face19_L.png'
$ws.Range('E107').Value = 'Face19_R.png'
$ws.Range('B108').Value = '3D'
$ws.Range('D108').Value = 'Face21_L.png'
$ws.Range('E108').Value = 'Face21_R.png'
$ws.Range('B109').Value = '3D'
$ws.Range('D109').Value = 'Face2_L.png'
$ws.Range('E109').Value = 'Face2_R.png'
$ws.Range('B110').Value = '3D'
$ws.Range('D110').Value = 'Face22_L.png'
$ws.Range('E110').Value = 'Face22_R.png'
$ws.Range('B111').Value = '3D'
$ws.Range('D111').Value = 'Face16_L.png'
$ws.Range('E111').Value = 'Face16_R.png'
$ws.Range('B112').Value = '3D'
$ws.Range('D112').Value = 'Face20_L.png'
$ws.Range('E112').Value = 'Face20_R.png'
$ws.Range('B113').Value = '3D'
$ws.Range('D113').Value = 'Face12_L.png'
$ws.Range('E113').Value = 'Face12_R.png'
$ws.Range('B114').Value = '3D'
$ws.Range('D114').Value = 'Face7_L.png'
$ws.Range('E114').Value = 'Face7_R.png'
$ws.Range('B115').Value = '3D'
$ws.Range('D115').Value = 'Face1_L.png'
$ws.Range('E115').Value = 'Face1_R.png'
$ws.Range('B116').Value = '3D'
$ws.Range('D116').Value = 'Face1_L.png'
$ws.Range('E116').Value = 'Face1_R.png'
$ws.Range('B117').Value = '3D'
$ws.Range('D117').Value = 'Face5_L.png'
$ws.Range('E117').Value = 'Face5_R.png'
$ws.Range('B118').Value = '3D'
$ws.Range('D118').Value = 'Face3_L.png'
$ws.Range('E118').Value = 'Face3_R.png'
$ws.Range('B119').Value = '3D'
$ws.Range('D119').Value = 'Face6_L.png'
$ws.Range('E119').Value = 'Face6_R.png'
$ws.Range('B120').Value = '3D'
$ws.Range('D120').Value = 'Face10_L.png'
$ws.Range('E120').Value = 'Face10_R.png'
$ws.Range('B121').Value = '3D'
$ws.Range('D121').Value = 'Face14_L.png'
$ws.Range('E121').Value = 'Face14_R.png'
$ws.Range('B122').Value = '3D'
$ws.Range('D122').Value = 'Face23_L.png'
$ws.Range('E122').Value = 'Face23_R.png'
$ws.Range('B123').Value = '3D'
$ws.Range('D123').Value = 'Face17_L.png'
$ws.Range('E123').Value = 'Face17_R.png'
$ws.Range('B124').Value = '3D'
$ws.Range('D124').Value = 'Face9_L.png'
$ws.Range('E124').Value = 'Face9_R.png'
$ws.Range('B125').Value = '3D'
$ws.Range('D125').Value = 'Face8_L.png'
$ws.Range('E125').Value = 'Face8_R.png'
$ws.Range('B126').Value = '3D'
$ws.Range('D126').Value = 'Face18_L.png'
$ws.Range('E126').Value = 'Face18_R.png'
$ws.Range('B127').Value = '3D'
$ws.Range('D127').Value = 'Face13_L.png'
$ws.Range('E127').Value = 'Face13_R.png'
$ws.Range('B128').Value = '3D'
$ws.Range('D128').Value = 'Face11_L.png'
$ws.Range('E128').Value = 'Face11_R.png'
$ws.Range('B129').Value = 'CUE'
$ws.Range('D129').Value = 'CoverCue.png'
$ws.Range('E129').Value = 'CoverCue.png'
$ws.Range('D131').Value = 'Face2_L.png'
$ws.Range('E131').Value = 'Face2_L.png'
$ws.Range('D132').Value = 'Face10_L.png'
$ws.Range('E132').Value = 'Face10_L.png'
$ws.Range('D133').Value = 'Face9_L.png'
$ws.Range('E133').Value = 'Face9_L.png'
$ws.Range('D134').Value = 'Face1_L.png'
$ws.Range('E134').Value = 'Face1_L.png'
$ws.Range('D135').Value = 'Face17_L.png'
$ws.Range('E135').Value = 'Face17_L.png'
$ws.Range('D136').Value = 'Face22_L.png'
$ws.Range('E136').Value = 'Face22_L.png'
$ws.Range('D137').Value = 'Face6_L.png'
$ws.Range('E137').Value = 'Face6_L.png'
$ws.Range('D138').Value = 'Face13_L.png'
$ws.Range('E138').Value = 'Face13_L.png'
$ws.Range('D139').Value = 'Face13_L.png'
$ws.Range('E139').Value = 'Face13_L.png'
$ws.Range('D140').Value = 'Face4_L.png'
$ws.Range('E140').Value = 'Face4_L.png'
$ws.Range('D141').Value = 'Face12_L.png'
$ws.Range('E141').Value = 'Face12_L.png'
$ws.Range('D142').Value = 'Face5_L.png'
$ws.Range('E142').Value = 'Face5_L.png'
$ws.Range('D143').Value = 'Face20_L.png'
$ws.Range('E143').Value = 'Face20_L.png'
$ws.Range('D144').Value = 'Face3_L.png'
$ws.Range('E144').Value = 'Face3_L.png'
$ws.Range('D145').Value = 'Face8_L.png'
$ws.Range('E145').Value = 'Face8_L.png'
$ws.Range('D146').Value = 'Face18_L.png'
$ws.Range('E146').Value = 'Face18_L.png'
$ws.Range('D147').Value = 'Face23_L.png'
$ws.Range('E147').Value = 'Face23_L.png'
$ws.Range('D148').Value = 'Face21_L.png'
$ws.Range('E148').Value = 'Face21_L.png'
$ws.Range('D149').Value = 'Face19_L.png'
$ws.Range('E149').Value = 'Face19_L.png'
$ws.Range('D150').Value = 'Face15_L.png'
$ws.Range('E150').Value = 'Face15_L.png'
$ws.Range('D151').Value = 'Face16_L.png'
$ws.Range('E151').Value = 'Face16_L.png'
$ws.Range('D152').Value = 'Face14_L.png'
$ws.Range('E152').Value = 'Face14_L.png'
$ws.Range('D153').Value = 'Face7_L.png'
$ws.Range('E153').Value = 'Face7_L.png'
$ws.Range('D154').Value = 'Face11_L.png'
$ws.Range('E154').Value = 'Face11_L.png'
$ws.Range('B155').Value = 'CUE'
$ws.Range('D155').Value = 'UncoverCue.png'
$ws.Range('E155').Value = 'UncoverCue.png'
$ws.Range('B157').Value = '2D'
$ws.Range('D157').Value = 'Face10_L.png'
$ws.Range('E157').Value = 'Face10_L.png'
$ws.Range('B158').Value = '2D'
$ws.Range('D158').Value = 'Face4_L.png'
$ws.Range('E158').Value = 'Face4_L.png'
$ws.Range('B159').Value = '2D'
$ws.Range('D159').Value = 'Face2_L.png'
$ws.Range('E159').Value = 'Face2_L.png'
$ws.Range('B160').Value = '2D'
$ws.Range('D160').Value = 'Face7_L.png'
$ws.Range('E160').Value = 'Face7_L.png'
$ws.Range('B161').Value = '2D'
$ws.Range('D161').Value = 'Face6_L.png'
$ws.Range('E161').Value = 'Face6_L.png'
$ws.Range('B162').Value = '2D'
$ws.Range('D162').Value = 'Face12_L.png'
$ws.Range('E162').Value = 'Face12_L.png'
$ws.Range('B163').Value = '2D'
$ws.Range('D163').Value = 'Face23_L.png'
$ws.Range('E163').Value = 'Face23_L.png'
$ws.Range('B164').Value = '2D'
$ws.Range('D164').Value = 'Face18_L.png'
$ws.Range('E164').Value = 'Face18_L.png'
$ws.Range('B165').Value = '2D'
$ws.Range('D165').Value = 'Face1_L.png'
$ws.Range('E165').Value = 'Face1_L.png'
$ws.Range('B166').Value = '2D'
$ws.Range('D166').Value = 'Face5_L.png'
$ws.Range('E166').Value = 'Face5_L.png'
$ws.Range('B167').Value = '2D'
$ws.Range('D167').Value = 'Face20_L.png'
$ws.Range('E167').Value = 'Face20_L.png'
$ws.Range('B168').Value = '2D'
$ws.Range('D168').Value = 'Face11_L.png'
$ws.Range('E168').Value = 'Face11_L.png'
$ws.Range('B169').Value = '2D'
$ws.Range('D169').Value = 'Face3_L.png'
$ws.Range('E169').Value = 'Face3_L.png'
$ws.Range('B170').Value = '2D'
$ws.Range('D170').Value = 'Face14_L.png'
$ws.Range('E170').Value = 'Face14_L.png'
$ws.Range('B171').Value = '2D'
$ws.Range('D171').Value = 'Face13_L.png'
$ws.Range('E171').Value = 'Face13_L.png'
$ws.Range('B172').Value = '2D'
$ws.Range('D172').Value = 'Face17_L.png'
$ws.Range('E172').Value = 'Face17_L.png'
$ws.Range('B173').Value = '2D'
$ws.Range('D173').Value = 'Face9_L.png'
$ws.Range('E173').Value = 'Face9_L.png'
$ws.Range('B174').Value = '2D'
$ws.Range('D174').Value = 'Face21_L.png'
$ws.Range('E174').Value = 'Face21_L.png'
$ws.Range('B175').Value = '2D'
$ws.Range('D175').Value = 'Face21_L.png'
$ws.Range('E175').Value = 'Face21_L.png'
$ws.Range('B176').Value = '2D'
$ws.Range('D176').Value = 'Face19_L.png'
$ws.Range('E176').Value = 'Face19_L.png'
$ws.Range('B177').Value = '2D'
$ws.Range('D177').Value = 'Face22_L.png'
$ws.Range('E177').Value = 'Face22_L.png'
$ws.Range('B178').Value = '2D'
$ws.Range('D178').Value = 'Face16_L.png'
$ws.Range('E178').Value = 'Face16_L.png'
$ws.Range('B179').Value = '2D'
$ws.Range('D179').Value = 'Face8_L.png'
$ws.Range('E179').Value = 'Face8_L.png'
$ws.Range('B180').Value = '2D'
$ws.Range('D180').Value = 'Face15_L.png'
$ws.Range('E180').Value = 'Face15_L.png'
$ws.Range('B182').Value = 'Pseudo'
$ws.Range('D182').Value = 'Face2_R.png'
$ws.Range('E182').Value = 'Face2_L.png'
$ws.Range('B183').Value = 'Pseudo'
$ws.Range('D183').Value = 'Face21_R.png'
$ws.Range('E183').Value = 'Face21_L.png'
$ws.Range('B184').Value = 'Pseudo'
$ws.Range('D184').Value = 'Face23_R.png'
$ws.Range('E184').Value = 'Face23_L.png'
$ws.Range('B185').Value = 'Pseudo'
$ws.Range('D185').Value = 'Face22_R.png'
$ws.Range('E185').Value = 'Face22_L.png'
$ws.Range('B186').Value = 'Pseudo'
$ws.Range('D186').Value = 'Face7_R.png'
$ws.Range('E186').Value = 'Face7_L.png'
$ws.Range('B187').Value = 'Pseudo'
$ws.Range('D187').Value = 'Face6_R.png'
$ws.Range('E187').Value = 'Face6_L.png'
$ws.Range('B188').Value = 'Pseudo'
$ws.Range('D188').Value = 'Face19_R.png'
$ws.Range('E188').Value = 'Face19_L.png'
$ws.Range('B189').Value = 'Pseudo'
$ws.Range('D189').Value = 'Face12_R.png'
$ws.Range('E189').Value = 'Face12_L.png'
$ws.Range('B190').Value = 'Pseudo'
$ws.Range('D190').Value = 'Face8_R.png'
$ws.Range('E190').Value = 'Face8_L.png'
$ws.Range('B191').Value = 'Pseudo'
$ws.Range('D191').Value = 'Face5_R.png'
$ws.Range('E191').Value = 'Face5_L.png'
$ws.Range('B192').Value = 'Pseudo'
$ws.Range('D192').Value = 'Face1_R.png'
$ws.Range('E192').Value = 'Face1_L.png'
$ws.Range('B193').Value = 'Pseudo'
$ws.Range('D193').Value = 'Face9_R.png'
$ws.Range('E193').Value = 'Face9_L.png'
$ws.Range('B194').Value = 'Pseudo'
$ws.Range('D194').Value = 'Face4_R.png'
$ws.Range('E194').Value = 'Face4_L.png'
$ws.Range('B195').Value = 'Pseudo'
$ws.Range('D195').Value = 'Face20_R.png'
$ws.Range('E195').Value = 'Face20_L.png'
$ws.Range('B196').Value = 'Pseudo'
$ws.Range('D196').Value = 'Face15_R.png'
$ws.Range('E196').Value = 'Face15_L.png'
$ws.Range('B197').Value = 'Pseudo'
$ws.Range('D197').Value = 'Face10_R.png'
$ws.Range('E197').Value = 'Face10_L.png'
$ws.Range('B198').Value = 'Pseudo'
$ws.Range('D198').Value = 'Face10_R.png'
$ws.Range('E198').Value = 'Face10_L.png'
$ws.Range('B199').Value = 'Pseudo'
$ws.Range('D199').Value = 'Face16_R.png'
$ws.Range('E199').Value = 'Face16_L.png'
$ws.Range('B200').Value = 'Pseudo'
$ws.Range('D200').Value = 'Face18_R.png'
$ws.Range('E200').Value = 'Face18_L.png'
$ws.Range('B201').Value = 'Pseudo'
$ws.Range('D201').Value = 'Face11_R.png'
$ws.Range('E201').Value = 'Face11_L.png'
$ws.Range('B202').Value = 'Pseudo'
$ws.Range('D202').Value = 'Face3_R.png'
$ws.Range('E202').Value = 'Face3_L.png'
$ws.Range('B203').Value = 'Pseudo'
$ws.Range('D203').Value = 'Face14_R.png'
$ws.Range('E203').Value = 'Face14_L.png'
$ws.Range('B204').Value = 'Pseudo'
$ws.Range('D204').Value = 'Face13_R.png'
$ws.Range('E204').Value = 'Face13_L.png'
$ws.Range('B205').Value = 'Pseudo'
$ws.Range('D205').Value = 'Face17_R.png'
$ws.Range('E205').Value = 'Face17_L.png'

# Reset the selection to span the full used data range
$ws.Range("A1:Q206").Select()

Write-Host "Applied image filename updates (jpg -> png) across the order sheet"